$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns Q ("Time") and R ("Power") ---
# Shared-string append order matters: "Power" must be inserted before "Time"
# so they land at sharedStrings indices 124 and 125 respectively.
$ws.Range("R1").Value = "Power"
$ws.Range("Q1").Value = "Time"

# --- Per-row Power (R) / Time (Q) numbers, rows 2..63 (30 and 47 are section header rows and stay empty) ---
$qValues = @(7,9,6,5,8,10,4,8,8,7,5,7,3,3,3,4,4,5,6,8,6,3,9,6,7,9,3,10,$null,4,5,3,4,5,4,6,4,6,6,4,7,6,8,7,8,$null,9,2,3,4,10,5,3,6,5,7,8,4,3,6,7,8)
$rValues = @(5,2,8,7,3,3,9,5,4,5,7,3,9,9,6,7,8,7,5,4,5,8,3,3,5,4,9,3,$null,9,8,9,8,7,8,7,6,6,5,8,4,3,2,2,4,$null,2,7,9,8,3,2,10,6,7,4,5,6,9,4,5,3)

for ($i = 0; $i -lt $qValues.Length; $i++) {
    $row = $i + 2
    $q = $qValues[$i]
    $r = $rValues[$i]
    if ($null -ne $q) {
        $ws.Cells.Item($row, 17).Value = $q
    }
    if ($null -ne $r) {
        $ws.Cells.Item($row, 18).Value = $r
    }
}

# --- Column widths / visibility ---
# Column O (15): new, visible
$ws.Columns.Item(15).ColumnWidth = 6.714285714285714
# Column P (16): new, hidden
$ws.Columns.Item(16).ColumnWidth = 8.142857142857142
$ws.Columns.Item(16).Hidden = $true
# Column Q (17): widened from 4.5 to ~11.2
$ws.Columns.Item(17).ColumnWidth = 10.428571428571429

# --- Sheet view: scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 42
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("Q63").Select()
